$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "fig" labels in column E for rows 31-35
$ws.Range("E31").Value = "12-fig-01"
$ws.Range("E32").Value = "13-fig-01"
$ws.Range("E33").Value = "13-fig-02"
$ws.Range("E34").Value = "13-fig-03"
$ws.Range("E35").Value = "13-fig-06"

# Scroll the view and move the selection, as captured in the saved sheetView
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("E36").Select()
